$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells that flip from their current value to the shared "NaN" text marker ---
$ws.Range("CU17").Value = "NaN"
$ws.Range("L18").Value = "NaN"
$ws.Range("BQ18").Value = "NaN"
$ws.Range("DB31").Value = "NaN"
$ws.Range("AK33").Value = "NaN"
$ws.Range("AK34").Value = "NaN"
$ws.Range("AK35").Value = "NaN"
$ws.Range("CJ36").Value = "NaN"
$ws.Range("BU45").Value = "NaN"
$ws.Range("AR57").Value = "NaN"
$ws.Range("AR58").Value = "NaN"
$ws.Range("CP75").Value = "NaN"
$ws.Range("DH93").Value = "NaN"
$ws.Range("DH94").Value = "NaN"
$ws.Range("AP102").Value = "NaN"
$ws.Range("BU106").Value = "NaN"
$ws.Range("AF147").Value = "NaN"

# --- Cells that flip from the "NaN" text marker back to a plain number ---
$ws.Range("BU9").Value = 1
$ws.Range("H90").Value = 1

# --- Append new row 168 with the day's report ---
$row168 = @{
    "A"  = 44062;  "B"  = 502178; "C"  = 2673;   "D"  = 65786;  "E"  = 62009
    "F"  = 176336; "G"  = 22485;  "H"  = 2230;   "I"  = 1825;   "J"  = 4339
    "K"  = 3512;   "L"  = 6439;   "M"  = 3547;   "N"  = 15786;  "O"  = 16756
    "P"  = 3891;   "Q"  = 2793;   "R"  = 10666;  "S"  = 4943;   "T"  = 11875
    "U"  = 7501;   "V"  = 2207;   "W"  = 681;    "X"  = 3718;   "Y"  = 11283
    "Z"  = 9440;   "AA" = 4970;   "AB" = 40073;  "AC" = 745;    "AD" = 112
    "AE" = 163;    "AF" = 435;    "AG" = 20;     "AH" = 14;     "AI" = 198
    "AJ" = 1917;   "AK" = 1973;   "AL" = 34965;  "AM" = 5440;   "AN" = 2330
    "AO" = 31545;  "AP" = 771;    "AQ" = 18565;  "AR" = 1374;   "AS" = 5228
    "AT" = 1346;   "AU" = 1525;   "AV" = 2793;   "AW" = 1282;   "AX" = 925
    "AY" = 2441;   "AZ" = 2544;   "BA" = 37674;  "BB" = 10258;  "BC" = 1553
    "BD" = 6213;   "BE" = 2396;   "BF" = 273;    "BG" = 1353;   "BH" = 2472
    "BI" = 723;    "BJ" = 1869;   "BK" = 7248;   "BL" = 6781;   "BM" = 6375
    "BN" = 13393;  "BO" = 1830;   "BP" = 727;    "BQ" = 4657;   "BR" = 4134
    "BS" = 4620;   "BT" = 960;    "BU" = 1200;   "BV" = 1770;   "BW" = 2098
    "BX" = 505;    "BY" = 3601;   "BZ" = 2023;   "CA" = 832;    "CB" = 567
    "CC" = 1556;   "CD" = 1615;   "CE" = 748;    "CF" = 670;    "CG" = 3650
    "CH" = 893;    "CI" = 971;    "CJ" = 969;    "CK" = 1266;   "CL" = 1084
    "CM" = 1011;   "CN" = 942;    "CO" = 912;    "CP" = 1006;   "CQ" = 459
    "CR" = 2727;   "CS" = 741;    "CT" = 707;    "CU" = 619;    "CV" = 1004
    "CW" = 893;    "CX" = 520;    "CY" = 635;    "CZ" = 667;    "DA" = 970
    "DB" = 803;    "DC" = 884;    "DD" = 710;    "DE" = 309;    "DF" = 318
    "DG" = 597;    "DH" = 466;    "DI" = 360;    "DJ" = 500;    "DK" = 289
    "DL" = 496;    "DM" = 686;    "DN" = 495;    "DO" = 466;    "DP" = 341
    "DQ" = 507;    "DR" = 105490; "DS" = 209749; "DT" = 6497;   "DU" = 91002
    "DV" = 59007;  "DW" = 18833;  "DX" = 6566
}

$order = "A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ","BA","BB","BC","BD","BE","BF","BG","BH","BI","BJ","BK","BL","BM","BN","BO","BP","BQ","BR","BS","BT","BU","BV","BW","BX","BY","BZ","CA","CB","CC","CD","CE","CF","CG","CH","CI","CJ","CK","CL","CM","CN","CO","CP","CQ","CR","CS","CT","CU","CV","CW","CX","CY","CZ","DA","DB","DC","DD","DE","DF","DG","DH","DI","DJ","DK","DL","DM","DN","DO","DP","DQ","DR","DS","DT","DU","DV","DW","DX"

foreach ($col in $order) {
    $ws.Range("$col" + "168").Value = $row168[$col]
}

# --- Update the visible selection to the new bottom-right corner of the sheet ---
[void]$ws.Range("DX168").Select()
